# Update "Provincias Spain" sheet:
#  - bump the "Datos actualizados" timestamp in A1
#  - the city/province ranking table (rows 22-32, 43-46, 54-57, 61-63) was
#    re-sorted and refreshed with new case counts; write the new
#    city name (col A) + Casos totales/activos/Recuperados/Muertes (cols B-E)
#    for every row whose data moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 16:20"

# row -> City, Casos totales, Casos activos, Recuperados, Muertes
$rows = @(
    @(22, "Cantabria",           1023, 25,   972, 26),
    @(23, "Gipuzkoa/Guipuzcoa",  1017, 1503, 630, 34),
    @(24, "Caceres",              991, 11,   896, 84),
    @(25, "Sevilla",              912, 17,   865, 30),
    @(26, "Aragon",               907, 29,   838, 40),
    @(27, "Valladolid",           886, 127,  702, 57),
    @(28, "Murcia",               872, 16,   836, 20),
    @(29, "Granada",              860, 13,   790, 57),
    @(30, "Leon",                 821, 118,  626, 77),
    @(31, "Burgos",               719, 137,  533, 49),
    @(32, "La Palma",             712, 30,   676, 2),
    @(43, "Fuerteventura",        288, 30,   267, 0),
    @(44, "Lugo",                 270, 153,  244, 4),
    @(45, "Palencia",             262, 28,   221, 13),
    @(46, "Cuenca",               253, 236,  180, 55),
    @(54, "Lanzarote",             49, 30,    44, 3),
    @(55, "Melilla",               48, 0,     47, 1),
    @(56, "La Gomera",             35, 30,    29, 0),
    @(57, "El Hierro",             30, 30,    30, 0),
    @(61, "Gran Canaria",           8, 30,     7, 11),
    @(62, "Arroyo de la Luz",       7, 0,      7, 0),
    @(63, "Tenerife",               3, 30,     3, 36)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
}
